# Updated documentation for MeanGermTime 2
# Rewords the "observation"-based formula explanations for MeanGermTime,
# VarGermTime and SEGermTime (Sheet1 rows 9-11, column C) to the
# "interval"-based wording, and moves the selection/active cell to C12
# (matching the author's edit position after updating row 12's docs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C9").Value = 'It is the average length of time required for maximum germination of a seed lot and is estimated according to the following formula.
$$\overline{T} = \frac{\sum_{i=1}^{k}N_{i}T_{i}}{\sum_{i=1}^{k}N_{i}}$$
Where, $T_{i}$ is the time from the start of the experiment to the $i$th interval, $N_{i}$ is the number of seeds germinated in the $i$th time interval (not the accumulated number, but the number corresponding to the $i$th interval) and $k$ is the total number of time intervals.
It is the inverse of mean germination rate ($\overline{V}$).
$$\overline{T} = \frac{1}{\overline{V}}$$'

$ws.Range("C10").Value = 'It is computed according to the following formula.
$$s_{T}^{2} = \frac{\sum_{i=1}^{k}N_{i}(T_{i}-\overline{T})^{2}}{\sum_{i=1}^{k}N_{i}-1}$$
Where, $T_{i}$ is the time from the start of the experiment to the $i$th interval, $N_{i}$ is the number of seeds germinated in the $i$th time interval (not the accumulated number, but the number corresponding to the $i$th interval) and $k$ is the total number of time intervals.'

$ws.Range("C11").Value = 'It signifies the accuracy of the calculation of the mean germination time.
It is estimated according to the following formula:
$$s_{\overline{T}} = \sqrt{\frac{s_{T}^{2}}{\sum_{i=1}^{k}N_{i}}}$$
Where, $N_{i}$ is the number of seeds germinated in the $i$th time interval (not the accumulated number, but the number corresponding to the $i$th interval) and $k$ is the total number of time intervals.'

$ws.Activate()
$ws.Range("C12").Select()
